# cryptos.xlsx refresh: price / 1h-volume update (+ a few ranking swaps)
# Matches the GitHub Actions bot commit "Updated cryptos list ... with GitHub Actions"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that must be stored as literal TEXT (prices in this
# sheet use a "thousands dot" look, e.g. "65.152.69", or are plain decimals
# like "2.30" that Excel would otherwise auto-coerce into a Number and lose
# the trailing zero). Temporarily forcing the Text number format during the
# write keeps the literal string, then resetting the cell style back to
# "Normal" drops the now-unneeded explicit format so the cell's style index
# is left exactly as it was before (no stray s="..." on the cell).
function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# --- Rows whose Price (D) and Volume 1h (E) changed, ranking unchanged ----
$updates = @(
    @{ Row = 2;  D = "65.152.69";   E = "  +0.73%  " },
    @{ Row = 3;  D = "3.577.09";    E = "  +4.54%  " },
    @{ Row = 4;  D = $null;         E = "  -0.14%  " },
    @{ Row = 5;  D = "599.63";      E = "  +3.32%  " },
    @{ Row = 6;  D = "138.66";      E = "  +3.81%  " },
    @{ Row = 7;  D = "3.573.31";    E = "  +4.61%  " },
    @{ Row = 8;  D = $null;         E = "  +0.08%  " },
    @{ Row = 9;  D = "0.496";       E = "  +3.04%  " },
    @{ Row = 10; D = "0.124";       E = "  +3.47%  " },
    @{ Row = 11; D = "7.01";        E = "  -0.63%  " },
    @{ Row = 12; D = "0.389";       E = "  +4.60%  " },
    @{ Row = 13; D = "4.185.71";    E = "  +4.39%  " },
    @{ Row = 14; D = "0.0000184";   E = "  +3.79%  " },
    @{ Row = 15; D = "27.52";       E = "  +5.68%  " },
    @{ Row = 16; D = "3.569.70";    E = "  +3.85%  " },
    @{ Row = 17; D = $null;         E = "  +1.27%  " },
    @{ Row = 18; D = "65.000.98";   E = "  +0.44%  " },
    @{ Row = 19; D = "10.23";       E = "  +8.69%  " },
    @{ Row = 20; D = "5.89";        E = "  +3.39%  " },
    @{ Row = 21; D = "14.44";       E = "  +7.82%  " },
    @{ Row = 22; D = "393.23";      E = "  +3.69%  " },
    @{ Row = 23; D = "0.581";       E = "  +8.08%  " },
    @{ Row = 24; D = "3.719.73";    E = "  +4.33%  " },
    @{ Row = 25; D = "74.22";       E = "  +3.63%  " },
    @{ Row = 26; D = $null;         E = "  +0.11%  " },
    @{ Row = 27; D = "0.0000117";   E = "  +13.43%  " },
    @{ Row = 28; D = "7.80";        E = "  +9.74%  " },
    @{ Row = 29; D = $null;         E = "  +0.45%  " },
    @{ Row = 30; D = "2.30";        E = "  +6.16%  " },
    @{ Row = 31; D = "8.36";        E = "  +5.59%  " },
    @{ Row = 32; D = "3.589.31";    E = "  +4.16%  " },
    @{ Row = 33; D = $null;         E = "  +24.29%  " },
    @{ Row = 34; D = "24.10";       E = "  +5.73%  " },
    @{ Row = 36; D = "0.145";       E = "  +2.33%  " },
    @{ Row = 39; D = "6.95";        E = "  +4.63%  " },
    @{ Row = 40; D = "5.09";        E = "  +12.68%  " },
    @{ Row = 41; D = "0.0815";      E = "  +8.04%  " },
    @{ Row = 44; D = "42.65";       E = "  +2.09%  " },
    @{ Row = 47; D = "4.49";        E = "  +5.80%  " },
    @{ Row = 48; D = "1.69";        E = "  +5.57%  " },
    @{ Row = 49; D = "2.494.90";    E = "  +14.11%  " },
    @{ Row = 50; D = "6.95";        E = "  +8.17%  " },
    @{ Row = 51; D = "303.82";      E = "  +11.16%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        Set-TextValue $ws.Cells.Item($u.Row, 4) $u.D
    }
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}

# --- Rows that swapped rank order (coin name / link / price / volume all
#     move together) ------------------------------------------------------
function Set-CoinRow {
    param($row, $name, $link, $price, $volume)
    $ws.Cells.Item($row, 2).Value = $name
    $ws.Cells.Item($row, 3).Value = $link
    Set-TextValue $ws.Cells.Item($row, 4) $price
    $ws.Cells.Item($row, 5).Value = $volume
}

Set-CoinRow 37 "ImmutableX" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx" "1.57" "  +9.03%  "
Set-CoinRow 38 "Monero" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr" "170.07" "  +0.61%  "

Set-CoinRow 42 "Mantle" "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt" "0.829" "  +3.60%  "
Set-CoinRow 43 "EnergySwap" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens" "26.77" "  +20.00%  "

Set-CoinRow 45 "ONDO" "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo" "1.24" "  +11.73%  "
Set-CoinRow 46 "FirstDigitalUSD" "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd" "0.999" "  -0.15%  "
